# "fixing facies analysis for publication (in progress)"
#
# Corrects the TOC (column O) values on Sheet1: several rows' TOC figures
# were recomputed, and a few samples (rows 2, 8, 9, 15, 17) no longer have
# a valid TOC measurement and are cleared back to blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Samples whose TOC value is not available anymore -> blank the cell
$ws.Range("O2").ClearContents()
$ws.Range("O8").ClearContents()
$ws.Range("O9").ClearContents()
$ws.Range("O15").ClearContents()
$ws.Range("O17").ClearContents()

# Samples with corrected TOC values
$ws.Range("O3").Value = 1.3
$ws.Range("O4").Value = 1.4
$ws.Range("O5").Value = 1.4
$ws.Range("O6").Value = 2.5
$ws.Range("O7").Value = 15.4
$ws.Range("O10").Value = 3.7
$ws.Range("O11").Value = 3.01
$ws.Range("O12").Value = 8.7
$ws.Range("O13").Value = 7.6
$ws.Range("O14").Value = 6.1
$ws.Range("O16").Value = 9.2
